{"js": "// The commit's XML diff (document.xml, footer1.xml, footnotes.xml, styles.xml)\n// is purely cosmetic: every changed line is the exact same element with the\n// exact same set of attributes/namespace declarations, just re-serialized in\n// (alphabetically) sorted order. No text, formatting, structure, or document\n// property actually changes between \"before\" and \"after\".\n//\n// There is therefore no user-visible / model-level edit to perform through\n// the Word JavaScript API: Office.js has no surface for controlling the\n// attribute/namespace emission order of the underlying OOXML, and doing so\n// would not correspond to any real document change anyway.\n//\n// We still touch the document body (load/sync) so the API is exercised, but\n// we intentionally make no modifications, keeping the document semantically\n// identical to the input - matching the (semantically empty) diff.\ncontext.document.body.load(\"text\");\nawait context.sync();\n", "ps1": "# The commit's XML diff (document.xml, footer1.xml, footnotes.xml, styles.xml)\n# is purely cosmetic: every changed line is the exact same element with the\n# exact same set of attributes/namespace declarations, just re-serialized in\n# (alphabetically) sorted order. No text, formatting, structure, or document\n# property actually changes between \"before\" and \"after\".\n#\n# There is therefore no user-visible / model-level edit to perform through the\n# Word COM object model: the object model has no surface for controlling the\n# attribute/namespace emission order of the underlying OOXML, and doing so\n# would not correspond to any real document change anyway.\n#\n# We still touch the document (read-only access) so the object model is\n# exercised, but we intentionally make no modifications, keeping the document\n# semantically identical to the input - matching the (semantically empty) diff.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
